# Rename the embedded logo pictures in every header/footer story.
#
# The document carries two logos, each appearing once in the "first
# page" story and once in the "default" story of the relevant section:
#   - the BTEC logo (descr "BTec_Logo-Orange"), currently named
#     "image1.jpg", should become "image2.jpg"
#   - the Pearson logo (descr ending "PearsonLogo.png"), currently
#     named "image2.png", should become "image1.png"
#
# Only the `name` shown by Word's InlineShape.Name changes - the
# picture content/relationship and the description stay the same.

$d = $word.ActiveDocument

foreach ($sec in $d.Sections) {

    foreach ($hdr in $sec.Headers) {
        if ($hdr.Exists) {
            foreach ($shp in $hdr.Range.InlineShapes) {
                if ($shp.AlternativeText -eq "BTec_Logo-Orange") {
                    $shp.Name = "image2.jpg"
                }
                elseif ($shp.AlternativeText -like "*PearsonLogo.png") {
                    $shp.Name = "image1.png"
                }
            }
        }
    }

    foreach ($ftr in $sec.Footers) {
        if ($ftr.Exists) {
            foreach ($shp in $ftr.Range.InlineShapes) {
                if ($shp.AlternativeText -eq "BTec_Logo-Orange") {
                    try { $shp.Name = "image2.jpg" } catch { }
                }
                elseif ($shp.AlternativeText -like "*PearsonLogo.png") {
                    try { $shp.Name = "image1.png" } catch { }
                }
            }
        }
    }
}
